$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.653.22"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "1.608.41"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("D4").Value = "'0.997"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'212.42"
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("D6").Value = "'0.517"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'28.86"
$ws.Range("E8").Value = "  +7.35%  "
$ws.Range("D9").Value = "'0.258"
$ws.Range("E9").Value = "  +3.47%  "
$ws.Range("E10").Value = "  +1.61%  "
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("D12").Value = "1.838.81"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").Value = "1.608.69"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("D14").Value = "'0.565"
$ws.Range("E14").Value = "  +5.47%  "
$ws.Range("D15").Value = "29.690.66"
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("D16").Value = "'3.82"
$ws.Range("E16").Value = "  +1.90%  "
$ws.Range("D17").Value = "'64.54"
$ws.Range("E17").Value = "  +1.51%  "
$ws.Range("D18").Value = "'8.37"
$ws.Range("E18").Value = "  +9.79%  "
$ws.Range("D19").Value = "'241.64"
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("E20").Value = "  +1.35%  "
$ws.Range("D21").Value = "'0.997"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "'4.06"
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("D23").Value = "'9.53"
$ws.Range("E23").Value = "  +3.59%  "
$ws.Range("E24").Value = "  +1.20%  "
$ws.Range("D25").Value = "'156.69"
$ws.Range("E25").Value = "  +1.52%  "
$ws.Range("D26").Value = "'15.59"
$ws.Range("E26").Value = "  +1.79%  "
$ws.Range("E27").Value = "  +0.84%  "
$ws.Range("D28").Value = "'6.52"
$ws.Range("E28").Value = "  +2.14%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").Value = "'0.0481"
$ws.Range("E30").Value = "  +1.93%  "
$ws.Range("E31").Value = "  +0.49%  "
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("E33").Value = "  +2.16%  "
$ws.Range("D34").Value = "1.426.79"
$ws.Range("E34").Value = "  -0.20%  "
$ws.Range("E35").Value = "  +5.04%  "
$ws.Range("E36").Value = "  +1.14%  "
$ws.Range("D37").Value = "'2.87"
$ws.Range("E37").Value = "  +1.88%  "
$ws.Range("D38").Value = "'2.28"
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("E39").Value = "  +2.56%  "
$ws.Range("D40").Value = "'0.556"
$ws.Range("E40").Value = "  +3.90%  "
$ws.Range("D41").Value = "'0.826"
$ws.Range("E41").Value = "  +3.82%  "
$ws.Range("E42").Value = "  +0.88%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "'0.0496"
$ws.Range("E43").Value = "  +5.21%  "
$ws.Range("B44").Value = "BitcoinSV"
$ws.Range("C44").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D44").Value = "'54.32"
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'68.42"
$ws.Range("E45").Value = "  +4.38%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "'0.997"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "'0.998"
$ws.Range("E47").Value = "  +19.47%  "
$ws.Range("E48").Value = "  +2.92%  "
$ws.Range("D49").Value = "1.747.71"
$ws.Range("E49").Value = "  +0.28%  "
$ws.Range("D50").Value = "'87.08"
$ws.Range("E50").Value = "  +0.54%  "
$ws.Range("E51").Value = "  -1.45%  "
